# "meal updated 1 tarikh" -- fill in the 1 Oct meal counts (row 4) for each
# member and the bazar charge collected on 2 Oct (row 5, column F). Every
# other changed cell in the sheet (U4, B18, F35, K35:U35, K36:T36, K37:T37,
# B38) is a formula that recalculates automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 ("Ovi", 1-Oct): meal counts per member, columns K (Antor) .. T (Tawhid)
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 2
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 2

# Row 5 ("Dhrubo", 2-Oct): bazar charge
$ws.Range("F5").Value = 303

# Match the author's on-screen selection/scroll state when the file was saved
$wb.Windows.Item(1).ScrollRow = 16
$ws.Range("N41").Select()
